$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.438.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9997'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3734'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.32'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3634'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.256'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08125'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9999'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.607'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001272'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.290'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.631.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06893'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.519'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23.436.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.114'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.414'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.332'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.281'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.812.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.829'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02814'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2522'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.07244'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.126'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08761'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.371'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7075'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6527'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.332'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9987'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.015'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07970'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.202'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.51%  '

